# Dialogue JSONs CheckUp Fulfillment w/ RulesEngine
# Adds two new intent rows (17 and 18) to the intent-definition sheet:
#   Row 17: en.user.query.disease.treatment
#   Row 18: en.user.confirm.symptom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: en.user.query.disease.treatment ---
$ws.Range("A17").Value = "en.user.query.disease.treatment"
$ws.Range("E17").Value = "PHASE-CHECK"
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = "bot-disease"
$ws.Range("I17").Value = "Dev?"
$ws.Range("K17").Value = "What are the treatments for /@copd:bot-disease/?*How do you treat /@asthma:bot-disease/?*How do I manage symptoms of /@cad:bot-disease/?*How do I deal with /@lung cancer:bot-disease/?*What are the ways to treat /@valve disease:bot-disease/?*What do people do with /@heart failure:bot-disease/?*What do people deal with /@arrhythmia:bot-disease/?*"

# --- Row 18: en.user.confirm.symptom ---
$ws.Range("A18").Value = "en.user.confirm.symptom"
$ws.Range("E18").Value = "PHASE-CHECK"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = "bot-affirm"
$ws.Range("I18").Value = "Dev?"
$ws.Range("K18").Value = "/@Yes:bot-affirm/*/@Sure:bot-affirm/*/@Yeah:bot-affirm/*/@Ok:bot-affirm/*/@Alright:bot-affirm/*/@Agree:bot-affirm/*/@No:bot-affirm/*/@Nope:bot-affirm/*/@Nah:bot-affirm/*/@Disagree:bot-affirm/*/@By no means:bot-affirm/*/@Absolutely not:bot-affirm/"

# Update the view so the new last row is in focus, matching the saved workbook state
$ws.Range("J18").Select()
